$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at S (column 19). Excel carries the adjacent
# formatting into the new column exactly like the author's diff shows
# (S2 picks up style 21 from Q2/R2, S3 picks up 13, S4 picks up 19,
# S5 picks up 16), and extends dimension / row "spans" automatically.
$ws.Columns.Item(19).Insert()

# New "2022" column header in the year row
$ws.Range("S3").Value = 2022

# Revised figures for the manufacturing GVA share row (row 4)
$ws.Range("P4").Value = 13.7
$ws.Range("Q4").Value = 13.1
$ws.Range("R4").Value = 11.8
$ws.Range("S4").Value = 13.6

# Revised figures for the GVA-per-capita row (row 5)
$ws.Range("P5").Value = 13.6
$ws.Range("Q5").Value = 12.5
$ws.Range("R5").Value = 13.5
$ws.Range("S5").Value = 20

# Park the cursor on the newly added column's blank header-row cell,
# matching the author's recorded selection after the edit
$ws.Range("S2").Select()
